$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: insert new "Russia" column between Japan (K) and Saudi Arabia (old L).
# Saudi Arabia / USA shift from L/M to M/N.
# Row 1
$ws.Range("B1").Value = "`$ bold('All')"
$ws.Range("C1").Value = "`$ bold('Europe')"
$ws.Range("D1").Value = "France"
$ws.Range("E1").Value = "Germany"
$ws.Range("F1").Value = "Italy"
$ws.Range("G1").Value = "Poland"
$ws.Range("H1").Value = "Spain"
$ws.Range("I1").Value = "United Kingdom"
$ws.Range("J1").Value = "Switzerland"
$ws.Range("K1").Value = "Japan"
$ws.Range("L1").Value = "Russia"
$ws.Range("M1").Value = "Saudi Arabia"
$ws.Range("N1").Value = "USA"

# Row 2
$ws.Range("A2").Value = "Minimum tax of 2% on billionaires'`nwealth, in voluntary countries"
$ws.Range("B2").Value = 0.809187541057965
$ws.Range("C2").Value = 0.843978206156072
$ws.Range("D2").Value = 0.848427692052067
$ws.Range("E2").Value = 0.824709898398252
$ws.Range("F2").Value = 0.873001447429718
$ws.Range("G2").Value = 0.803221566293534
$ws.Range("H2").Value = 0.822122794211219
$ws.Range("I2").Value = 0.857393072880833
$ws.Range("J2").Value = 0.788973371785322
$ws.Range("K2").Value = 0.806207031821324
$ws.Range("L2").Value = 0.801008544697843
$ws.Range("M2").Value = 0.859553623962162
$ws.Range("N2").Value = 0.757626779099893

# Row 3
$ws.Range("A3").Value = "Bridgetown initiative: MDBs expanding sustainable`ninvestments in LICs, and at lower interest rates"
$ws.Range("B3").Value = 0.793366965036613
$ws.Range("C3").Value = 0.820231998216116
$ws.Range("D3").Value = 0.793373367472808
$ws.Range("E3").Value = 0.799097072664609
$ws.Range("F3").Value = 0.862342095281182
$ws.Range("G3").Value = 0.710759087360456
$ws.Range("H3").Value = 0.797966196862339
$ws.Range("I3").Value = 0.839503169105142
$ws.Range("J3").Value = 0.74891156552925
$ws.Range("K3").Value = 0.803594351085608
$ws.Range("L3").Value = 0.829224529123518
$ws.Range("M3").Value = 0.870851551659494
$ws.Range("N3").Value = 0.725878860090792

# Row 4
$ws.Range("A4").Value = "L&D: Developed countries financing a fund to help`nvulnerable countries cope with climate Loss and damage"
$ws.Range("B4").Value = 0.748295871658231
$ws.Range("C4").Value = 0.748627231605022
$ws.Range("D4").Value = 0.703301082880099
$ws.Range("E4").Value = 0.720152388692165
$ws.Range("F4").Value = 0.807893199188063
$ws.Range("G4").Value = 0.713408809337982
$ws.Range("H4").Value = 0.775759859493595
$ws.Range("I4").Value = 0.720550561248712
$ws.Range("J4").Value = 0.655019621133529
$ws.Range("K4").Value = 0.728151426836741
$ws.Range("L4").Value = 0.866768469161838
$ws.Range("M4").Value = 0.894499688071821
$ws.Range("N4").Value = 0.683507021804741

# Row 5
$ws.Range("A5").Value = "International levy on shipping carbon emissions,`nreturned to countries based on population"
$ws.Range("B5").Value = 0.699644803977167
$ws.Range("C5").Value = 0.732870289440397
$ws.Range("D5").Value = 0.752044185944775
$ws.Range("E5").Value = 0.6940206070794
$ws.Range("F5").Value = 0.754730850985435
$ws.Range("G5").Value = 0.606359310516872
$ws.Range("H5").Value = 0.747577003038328
$ws.Range("I5").Value = 0.738354993609978
$ws.Range("J5").Value = 0.708897221588476
$ws.Range("K5").Value = 0.576414783896346
$ws.Range("L5").Value = 0.72514029245591
$ws.Range("M5").Value = 0.814900578705803
$ws.Range("N5").Value = 0.650955072276642

# Row 6
$ws.Range("A6").Value = "At least 0.7% of developed countries' GDP in foreign aid"
$ws.Range("B6").Value = 0.698715666285492
$ws.Range("C6").Value = 0.689223510024568
$ws.Range("D6").Value = 0.640174218758149
$ws.Range("E6").Value = 0.662665049465059
$ws.Range("F6").Value = 0.765762693878515
$ws.Range("G6").Value = 0.589760130052602
$ws.Range("H6").Value = 0.76727730886302
$ws.Range("I6").Value = 0.684062495016513
$ws.Range("J6").Value = 0.648392506898434
$ws.Range("K6").Value = 0.611045583899474
$ws.Range("L6").Value = 0.825522087700358
$ws.Range("M6").Value = 0.863575793802146
$ws.Range("N6").Value = 0.640485413082898

# Row 7
$ws.Range("A7").Value = "Debt relief for vulnerable countries, suspending`npayments until they are more able to repay"
$ws.Range("B7").Value = 0.696886644817168
$ws.Range("C7").Value = 0.698508510826339
$ws.Range("D7").Value = 0.639494227176162
$ws.Range("E7").Value = 0.584238476063614
$ws.Range("F7").Value = 0.800152659041669
$ws.Range("G7").Value = 0.794100548139283
$ws.Range("H7").Value = 0.724925850622296
$ws.Range("I7").Value = 0.737663233832229
$ws.Range("J7").Value = 0.640198216978104
$ws.Range("K7").Value = 0.681263701142876
$ws.Range("L7").Value = 0.746898541844184
$ws.Range("M7").Value = 0.883360143018035
$ws.Range("N7").Value = 0.658719703625068

# Row 8
$ws.Range("A8").Value = "Expand Security Council to new permanent members (e.g.`nIndia, Brazil, African Union), restrict veto use"
$ws.Range("B8").Value = 0.694962460493881
$ws.Range("C8").Value = 0.761932435910776
$ws.Range("D8").Value = 0.721144190433339
$ws.Range("E8").Value = 0.753374080014139
$ws.Range("F8").Value = 0.760760292004357
$ws.Range("G8").Value = 0.704753173076338
$ws.Range("H8").Value = 0.761826274893759
$ws.Range("I8").Value = 0.795174607333805
$ws.Range("J8").Value = 0.70099140405434
$ws.Range("K8").Value = 0.674878873455674
$ws.Range("L8").Value = 0.52836576507021
$ws.Range("M8").Value = 0.836791104476278
$ws.Range("N8").Value = 0.649041911884122

# Row 9
$ws.Range("A9").Value = "NCQG: Developing countries providing `$300 bn a`nyear in climate finance for developing countries"
$ws.Range("B9").Value = 0.683174203642518
$ws.Range("C9").Value = 0.69192411456902
$ws.Range("D9").Value = 0.67136180939665
$ws.Range("E9").Value = 0.679657553482977
$ws.Range("F9").Value = 0.719275169766716
$ws.Range("G9").Value = 0.636020686294837
$ws.Range("H9").Value = 0.731089662519959
$ws.Range("I9").Value = 0.67838218523904
$ws.Range("J9").Value = 0.641704396059853
$ws.Range("K9").Value = 0.579399788140582
$ws.Range("L9").Value = 0.875589503017251
$ws.Range("M9").Value = 0.858825470219245
$ws.Range("N9").Value = 0.591520840537896

# Row 10
$ws.Range("A10").Value = "Raise global minimum tax on profit from 15% to 35%,`nallocating revenues to countries based on sales"
$ws.Range("B10").Value = 0.682119272063001
$ws.Range("C10").Value = 0.746931207916325
$ws.Range("D10").Value = 0.728051220092119
$ws.Range("E10").Value = 0.730873843089586
$ws.Range("F10").Value = 0.845480597940989
$ws.Range("G10").Value = 0.665987036677
$ws.Range("H10").Value = 0.685272115193463
$ws.Range("I10").Value = 0.736136988028249
$ws.Range("J10").Value = 0.634129522997038
$ws.Range("K10").Value = 0.730988448179478
$ws.Range("L10").Value = 0.50012458723726
$ws.Range("M10").Value = 0.770097780765828
$ws.Range("N10").Value = 0.631762229275282

# Row 11
$ws.Range("A11").Value = "International levy on aviation carbon emissions, raising`nprices by 30%, returned to countries based on population"
$ws.Range("B11").Value = 0.526319576644156
$ws.Range("C11").Value = 0.551475483086733
$ws.Range("D11").Value = 0.613989879120546
$ws.Range("E11").Value = 0.538994299398148
$ws.Range("F11").Value = 0.52078830506763
$ws.Range("G11").Value = 0.478870068276658
$ws.Range("H11").Value = 0.50585155387453
$ws.Range("I11").Value = 0.513338434920881
$ws.Range("J11").Value = 0.501802322471987
$ws.Range("K11").Value = 0.457604805676769
$ws.Range("L11").Value = 0.514246894490421
$ws.Range("M11").Value = 0.695189510992556
$ws.Range("N11").Value = 0.475985288378926
